$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for a new participant ("Ajay") above the current row 3
# (Deepa), shifting Deepa and NewUser down by one row.
$ws.Range("A3:H3").EntireRow.Insert()

# Populate the newly inserted row with the participant's data.
$ws.Cells.Item(3, 1).Value = 5.0
$ws.Cells.Item(3, 2).Value = "Ajay"
$ws.Cells.Item(3, 3).Value = "I05235"
$ws.Cells.Item(3, 4).Value = 100.0
$ws.Cells.Item(3, 5).Value = 5.0
$ws.Cells.Item(3, 6).Value = 5.0
$ws.Cells.Item(3, 7).Value = 26.0
$ws.Cells.Item(3, 8).Value = "2025-04-27 13:06:43"
